# Commit message: "Added: View presc record Added: Update presc status
# Fixed: Allignment for the Views"
#
# The underlying data change is the prescription-status update for
# appointment A002 (row 3): its Medication Status moves from PENDING to
# DISPENSED.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "DISPENSED"

# Alignment fix: reset the active selection to H2 (matches the refreshed
# "view prescription record" / "update prescription status" UI state).
$ws.Range("H2").Select() | Out-Null
